$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 289 ("Fruta / hortaliza, semanal"),
# shifting every subsequent row (old 289..379) down by one (new 290..380).
$ws.Rows.Item(289).Insert()

$ws.Range("A289").Value = 5
$ws.Range("B289").Value = "Macroferia Regional de Talca"
$ws.Range("C289").Value = "Maule"
$ws.Range("D289").Value = 44876
$ws.Range("E289").Value = 7
$ws.Range("F289").Value = 100112003
$ws.Range("G289").Value = "Ajo"
$ws.Range("H289").Value = "Chino"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 200
$ws.Range("K289").Value = 18000
$ws.Range("L289").Value = 18000
$ws.Range("M289").Value = 18000
$ws.Range("N289").Value = "$/malla 10 kilos"
$ws.Range("O289").Value = "China"
$ws.Range("P289").Value = 1800
$ws.Range("Q289").Value = 10
$ws.Range("R289").Value = "Hortaliza"
